$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'44.430.47"
$ws.Range('E2').Value = '  +0.91%  '
$ws.Range('D3').Value = "'2.234.31"
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  +1.66%  '
$ws.Range('D5').Value = "'306.29"
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('D6').Value = "'93.45"
$ws.Range('E6').Value = '  -2.19%  '
$ws.Range('D7').Value = "'0.570"
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('D9').Value = "'0.519"
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('D10').Value = "'34.45"
$ws.Range('E10').Value = '  -1.37%  '
$ws.Range('D11').Value = "'0.0798"
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('D12').Value = "'7.16"
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('D14').Value = "'2.265.11"
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').Value = "'0.828"
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('D16').Value = "'13.45"
$ws.Range('E16').Value = '  -0.72%  '
$ws.Range('D17').Value = "'44.091.52"
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('D18').Value = "'" + '0.0' + [char]0x2083 + '0949'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').Value = "'6.31"
$ws.Range('E19').Value = '  +1.49%  '
$ws.Range('D20').Value = "'11.88"
$ws.Range('E20').Value = '  -3.16%  '
$ws.Range('D21').Value = "'65.66"
$ws.Range('E21').Value = '  +1.35%  '
$ws.Range('D22').Value = "'236.72"
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').Value = "'2.95"
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('D24').Value = "'1.96"
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('D26').Value = "'2.22"
$ws.Range('E26').Value = '  +3.77%  '
$ws.Range('D27').Value = "'9.72"
$ws.Range('E27').Value = '  -2.38%  '
$ws.Range('D28').Value = "'37.40"
$ws.Range('E28').Value = '  -1.86%  '
$ws.Range('D29').Value = "'5.88"
$ws.Range('E29').Value = '  -1.16%  '
$ws.Range('D30').Value = "'19.83"
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('D31').Value = "'153.22"
$ws.Range('E31').Value = '  -1.19%  '
$ws.Range('D32').Value = "'0.0790"
$ws.Range('E32').Value = '  -2.16%  '
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('D34').Value = "'3.08"
$ws.Range('E34').Value = '  -7.02%  '
$ws.Range('E35').Value = '  +1.07%  '
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('D37').Value = "'1.77"
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').Value = "'3.36"
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('D39').Value = "'14.35"
$ws.Range('E39').Value = '  -5.88%  '
$ws.Range('D40').Value = "'3.74"
$ws.Range('E40').Value = '  -2.05%  '
$ws.Range('D41').Value = "'0.0298"
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('E42').Value = '  +0.74%  '
$ws.Range('D43').Value = "'1.761.10"
$ws.Range('E43').Value = '  +1.42%  '
$ws.Range('D44').Value = "'0.191"
$ws.Range('E44').Value = '  +1.53%  '
$ws.Range('D45').Value = "'78.62"
$ws.Range('E45').Value = '  -7.87%  '
$ws.Range('B46').Value = 'THORChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D46').Value = "'4.88"
$ws.Range('E46').Value = '  -0.97%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = "'98.15"
$ws.Range('E47').Value = '  -1.80%  '
$ws.Range('D48').Value = "'69.21"
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = "'8.03"
$ws.Range('E49').Value = '  -0.55%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').Value = "'54.45"
$ws.Range('E50').Value = '  +0.39%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = "'1.56"
$ws.Range('E51').Value = '  +2.98%  '
